$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 updates
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0.001
$ws.Range("K4").Value = 471
$ws.Range("L4").Value = 0.000942

# Row 5 updates
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 0.001
$ws.Range("K5").Value = 472
$ws.Range("L5").Value = 0.000944
